$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell C1 = "percent", styled like the existing header cells (A1/B1):
# bold font, thin border on all sides, centered horizontally, top-aligned vertically.
$ws.Range("C1").Value = "percent"
$ws.Range("C1").Font.Bold = $true
$ws.Range("C1").HorizontalAlignment = -4108
$ws.Range("C1").VerticalAlignment = -4160
$ws.Range("C1").Borders.LineStyle = 1

# Update the existing DOI count in B2
$ws.Range("B2").Value = 4111

# New "percent" text values in column C, stored as text (not numeric)
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "41.13"

$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "58.87"
